$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unprotect the sheet so we can modify cells/values as needed
$ws.Unprotect()

# Update the confidentiality/date notice text
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-27 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) values for rows 2-13
$ws.Range("D2").Value = 0.03121548328701642
$ws.Range("E2").Value = -0.003024193548387233

$ws.Range("D3").Value = 0.0233697740252822
$ws.Range("E3").Value = -0.006399620763213987

$ws.Range("D4").Value = 0.05206375471124405
$ws.Range("E4").Value = -0.005310551835603805

$ws.Range("D5").Value = 0.137188951823331
$ws.Range("E5").Value = 0.0006494560805325733

$ws.Range("D6").Value = 0.03116292006383014
$ws.Range("E6").Value = 0.0007072135785006761

$ws.Range("D7").Value = 0.1177198168065603
$ws.Range("E7").Value = 0.004450095359186168

$ws.Range("D8").Value = 0.1017081868987934
$ws.Range("E8").Value = 0.01299890150128169

$ws.Range("D9").Value = 0.02914536441870497
$ws.Range("E9").Value = 0.009898948236749749

$ws.Range("D10").Value = 0.1260338808865507
$ws.Range("E10").Value = 0.01159528707686563

$ws.Range("D11").Value = 0.2465627659098034
$ws.Range("E11").Value = -0.003905185723367621

$ws.Range("D12").Value = 0.1038291011688835
$ws.Range("E12").Value = 0.005148741418764313

$ws.Range("E13").Value = 0.002758272814913809

# Re-protect the sheet to restore its original protection state
$ws.Protect("D382")
